$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove "Hawaii" from the nhpi_country list (cell B9)
$ws.Range("B9").Select()
$ws.Range("B9").ClearContents()
